$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a cell while forcing it to be stored as literal
# text (Excel would otherwise auto-convert strings like "75.9%" into a
# percentage number). We stage the text in a scratch cell formatted as
# Text, copy it, paste-special (values) into the destination (this keeps
# the destination's existing style/number-format untouched), then clear
# the scratch cell so it leaves no trace in the sheet.
function Set-TextValue {
    param($ws, $cellRef, $text)
    $tmp = $ws.Range("Z1")
    $tmp.NumberFormat = "@"
    $tmp.Value = $text
    $tmp.Copy()
    $ws.Range($cellRef).PasteSpecial(-4163)
    $tmp.Clear()
}

# --- Row 2: reorder "Recorded By" list ---
$ws.Range("G2").Value = "System, servinaz@med.asu.edu.eg, gehanadel@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg"

# --- Row 3: reorder "Recorded By" list ---
$ws.Range("G3").Value = "majorelle.magdy@med.asu.edu.eg, System, asmaa.reda@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg"

# --- Row 4: reorder "Recorded By" list ---
$ws.Range("G4").Value = "majorelle.magdy@med.asu.edu.eg, servinaz@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, gehanadel@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg"

# --- Row 5: reorder "Recorded By" list ---
$ws.Range("G5").Value = "Amira.Sobhy@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg"

# --- Row 6: reorder "Recorded By" list, and update Recorded Sessions metric ---
$ws.Range("G6").Value = "alshimaa.atef@med.asu.edu.egm, majorelle.magdy@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg"
$ws.Range("L6").Value = 22

# --- Row 7: reorder "Recorded By" list, and update Missing Sessions metric ---
$ws.Range("G7").Value = "AbeerRagheb@med.asu.edu.eg, lamiaa.ossama@med.asu.edu.eg, NadaMohamed@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg, Fatmaelhady@med.asu.edu.eg, Kerelos.zareef@med.asu.edu.eg, Amera.a.saad@med.asu.edu.eg"
$ws.Range("L7").Value = 2

# --- Row 9: update Coverage % metric ---
Set-TextValue $ws "L9" "75.9%"

# --- Row 10: update Average Attendance % metric ---
Set-TextValue $ws "L10" "27.3%"

# --- Row 12: reorder "Recorded By" list ---
$ws.Range("G12").Value = "Eman.m.abosakaya@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, dina.adel@med.asu.edu.eg"

# --- Row 15: reorder "Recorded By" list, and update group-statistics row ---
$ws.Range("G15").Value = "Rania.a.youssef@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg"
$ws.Range("O15").Value = 22
$ws.Range("P15").Value = 2
Set-TextValue $ws "R15" "75.9%"
Set-TextValue $ws "S15" "27.3%"

# --- Row 17: reorder "Recorded By" list ---
$ws.Range("G17").Value = "esraa.sami@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg"

# --- Row 24: session now recorded - copy the formatting of an already-"Recorded" row (23)
#     then fill in the recorded-by / students / status values ---
$ws.Range("A23:I23").Copy()
$ws.Range("A24:I24").PasteSpecial(-4122)
$ws.Range("G24").Value = "Sarah.Mahdy@med.asu.edu.eg, youstina.gamil@med.asu.edu.eg"
$ws.Range("H24").Value = "106/251"
$ws.Range("I24").Value = "Recorded"

# --- Row 27: reorder "Recorded By" list ---
$ws.Range("G27").Value = "hana.amr@med.asu.edu.eg, nourhan.mostafa@med.asu.edu.eg"

# --- Row 30: reorder "Recorded By" list ---
$ws.Range("G30").Value = "shorokmohamed@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg"
